# This script reorders the weekly price-report rows (rows 2-9) so that the
# Fecha/Volumen/Precio data moves to the correct week, per the commit
# "Fruta / hortaliza, semanal". Column A,B,C,E,F,G,H,I,J,K,L,R,T are
# unaffected; only D (Fecha), M (Volumen), N (Precio minimo),
# O (Precio maximo), P (Precio promedio ponderado),
# Q (Unidad de comercializacion) and S (Precio $/Kg) change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final values to write per row, taken directly from the target OOXML.
$rows = @{
    2 = @{ D = 44351; M = 300; N = 10000; O = 10000; P = 10000; Q = "`$/caja 14 kilos empedrada"; S = 714 }
    3 = @{ D = 44208; M = 210; N = 10000; O = 10000; P = 10000; Q = "`$/caja 14 kilos empedrada"; S = 714 }
    4 = @{ D = 44491; M = 180; N = 9000;  O = 9000;  P = 9000;  Q = "`$/caja 14 kilos empedrada"; S = 643 }
    5 = @{ D = 44309; M = 300; N = 7000;  O = 7000;  P = 7000;  Q = "`$/caja 14 kilos empedrada"; S = 500 }
    6 = @{ D = 44400; M = 100; N = 10000; O = 10000; P = 10000; Q = "`$/caja 14 kilos";           S = 714 }
    7 = @{ D = 44162; M = 120; N = 7000;  O = 7000;  P = 7000;  Q = "`$/caja 14 kilos empedrada"; S = 500 }
    8 = @{ D = 44397; M = 60;  N = 11000; O = 11000; P = 11000; Q = "`$/caja 14 kilos";           S = 786 }
    9 = @{ D = 44176; M = 250; N = 7000;  O = 7000;  P = 7000;  Q = "`$/caja 14 kilos empedrada"; S = 500 }
}

foreach ($r in $rows.Keys) {
    $v = $rows[$r]
    $ws.Range("D$r").Value = $v.D
    $ws.Range("M$r").Value = $v.M
    $ws.Range("N$r").Value = $v.N
    $ws.Range("O$r").Value = $v.O
    $ws.Range("P$r").Value = $v.P
    $ws.Range("Q$r").Value = $v.Q
    $ws.Range("S$r").Value = $v.S
}
